$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.993.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.642.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.84%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4764'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.69%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2601'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.51%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06122'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07036'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.649.26'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.34%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.58'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5898'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.86%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.343'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.32%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '73.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.09%  '
$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.991.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.94%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006602'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.98%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.857.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.78%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.283'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.65%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.560'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.241'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.09%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '133.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.386'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.55%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '103.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.635'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.76%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.888'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07667'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.19%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.568'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.06%  '
$ws.Range('B33').Value = 'Frax'
$ws.Range('C33').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04286'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.68%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.571'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5905'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.54%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9258'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.00%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.577'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.43%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8861'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.49%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01508'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.16%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.756'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.12%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3691'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.87%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.675'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.74%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1101'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.20%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.085'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.91%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05208'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '28.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.06%  '
$ws.Range('B50').Value = 'TrueUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9984'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
